# Week 13 logging: append this week's game-log entries to the running
# per-game stat strings (YDS / ST sheets) and roll the new totals into
# the season-total cells on OFF, DEF, ST, TURNS, and PEN.
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# YDS sheet - append Week 13 rushing/passing yard entries
# ---------------------------------------------------------------
$ydsWs = $wb.Worksheets.Item("YDS")
$ydsWs.Range("B2").Value = "7 10 5 4 2 5 4 2 15 -1 7 7 4 3 1 4 5 -2 6 5 6 7 11 1 8 5 2 11 2 2 23 7 4 12 5 -1 5 1 -1 11 3 3 5 1 1 2 5 1 1 8 2 1 4 3 0 -1 6 2 0 -4 13 1 12 -1 4 0 2 3 -1 2 9 3 5 13 6 5 4 17 5 8 13 3 -1 1 7 38 6 4 5 48 -1 0 2 2 21 3 14 2 2 8 5 2 2 3 4 10 18 -1 2 1 2 3 6 2 5 6 -1 8 4 4 7 0 9 41 4 2 5 3 5 1 10 6 2 12 7 1 6 0 1 -1 2 2 0 4 0 7 2 2 38 1 6 1 12 3 5 2 4 5 6 2 2 3 7 12 2 4 3 3 9 10 6 -3 2 2 2 6 7 4 9 3 -2 2 2 5 2 2 4 2 6 -2 2 3 6 8 5 10 5 22 13 18 7 1 2 6 2 19 3 6 1 5 12 3 0 5 1 1 3 5 10 12 -1 1 3 11 2 5 5 3 7 4 6 1 3 7 2 2 7 6 1 -4 2 3 21 2 0 1 4 5 7 4 7 3 0 6 7 4 3 4 4 12 5 1 13 4 16 25 4 12 3 4 1 -2 -1 0 9 1 1 0 -1 8 8 1 2 -2 3 4 7 2 1 9 7 4 11 3 4 5 3 3 -1 0 8 -2 -3 9 2 2 6 0 4 4 -2 0 5 6 10 1 8 4 3 -1 15 3 0 18 6 6 0 0 1 4 -8 1 2 0 14 2 7 8 11 14 1 4 4 1 4 1 14 10 2 0 3 12 3 1 5 1 4 2 4 1 2 3 6 3 2 7 5 3 3 4 -4 0 3 5 0 1 3 12 1 6 5 9 8 4 4 3 1 2 1 0 0 7 33 7 4 9 3 3 3 1 6 -1 -1 6 16 1 4 17 4 3 8 -1 4 2 3 4 2 4 11 10 2 13 2 3 3 29 -3 3 2 8 6 1 21 28 3 9 2 7 3 7 1 8 -2 2 2 3 2 0 13 49 3 10 3 -1 0 -1 4 2 5 6 3 5 4 5 5 9 2 8 2 7 12 1 9 1 2 -1 35 3 6 1 2 2 11 4 5 2 0 3 4 4 1 4 5 1 -1 10 -4 1 3 3 1 8 1 3 5 2 4 1 -1 9 2 5 2 7 7 2 4 2 26 9 16 1 -1 1 5 -2 12 3 -1 -5 7 -2 4 3 6 1 2 0 1 17 17 2 2 2 0 -1 7 12 -4 0 0 0 -1 1 4 4 12 0 7 0 1 7 4 -1 4 0 -3 3 -3 4 7 5 15 4 0 8 3 1 24 9 3 2 5 1 21 4 4 5 3 1 1 1 11 -1 6 9 1 3 4 1 21 6 6 3 1 0 0 6 5 2 12 5 9 2 32 6 1 3 10 10 0 13 3 0 -2 15 5 2 4 9 6 3 -2 1 6 -1 -2 0 1 3 -1 3 4 1 8 5 6 5 6 3 7 4 7 5 2 6 10 4 6 5 4 -2 6 3 2 14 -1 6 3 2 0 1 0 1 1 2 2 2 0 16 12 1 5 2 -1 6 2 13 3 3 2 8 4 11 12 2 -2 3 6 5 1 4 11 2 5 1 7 -6 1 2 9 0 5 -1 15 9 13 4 0 5 5 4 18 2 3 1 5 1 -1 16 11 5 1 5 22 17 2 5 10 -1 2 3 0 1 1 1 21 3 2 1 2 14 17 2 1 5 13 -1 7 -1 5 12 2 10 0 4 1 6 2 7 0 3 1 9 0 2 2 2 2 2 1 6 5 7 4 2 -1 4 6 10 19 -1 7 14 2 4 1 2 4 -1 -1 6 64 10 6 0 16 3 3 2 -2 -4 10 17 2 1 5 4 7 3 4 1 22 11 6 -1 4 1 10 5 5 -3 -5 3 -1 2 2"
$ydsWs.Range("C2").Value = "4 2 2 1 -3 7 3 5 6 0 0 8 10 2 5 4 -1 1 1 1 12 5 4 4 7 5 6 2 21 3 5 0 1 9 2 0 5 10 3 5 13 7 3 6 8 2 6 3 3 5 6 3 12 6 7 6 6 5 -1 2 23 13 3 1 2 4 4 3 2 -1 9 3 14 3 6 5 1 3 13 1 9 4 2 1 5 -1 11 4 1 2 1 18 9 4 2 4 3 4 4 2 11 9 2 3 20 2 9 3 6 4 0 5 1 9 2 1 4 1 -2 3 -2 8 0 3 -3 1 9 14 -1 0 3 4 4 2 3 -2 5 3 13 9 1 4 7 17 6 3 4 20 5 16 8 5 7 3 6 6 9 10 -1 1 2 0 9 3 2 1 -1 3 6 3 0 19 21 8 -2 2 4 6 1 2 7 6 6 18 4 2 2 5 9 17 2 12 3 -2 2 2 3 8 14 1 -1 5 4 3 2 3 3 2 15 3 7 1 2 3 3 3 1 0 5 -1 4 4 3 1 5 7 5 3 3 9 5 1 3 5 6 1 2 1 8 5 -2 7 11 0 8 2 1 -2 1 2 8 -1 7 4 11 3 3 3 3 2 3 3 1 1 3 14 7 2 -2 4 2 0 3 3 3 4 2 9 -1 9 5 2 4 4 15 4 13 6 2 1 1 5 4 4 4 1 9 2 7 -1 1 3 3 6 8 5 8 5 6 6 35 2 4 14 9 6 11 2 1 1 4 6 6 8 25 3 2 1 9 1 4 0 2 3 3 0 3 6 3 4 -2 -1 1 5 6 5 4 6 7 2 10 4 8 4 6 3 4 5 4 9 14 2 7 5 1 31 8 3 2 0 12 3 4 3 13 24 9 4 2 1 0 1 2 6 3 2 4 -1 5 3 11 5 2 3 11 4 4 1 3 3 2 2 11 6 2 5 4 2 3 12 4 9 3 4 7 3 3 3 1 8 10 1 0 1 4 15 1 5 2 8 5 5 4 2 3 2 15 4 -1 8 3 -1 -1 2 2 1 15 8 4 -3 5 2 5 4 2 8 13 3 6 1 3 7 3 3 -1 17 5 1 4 4 6 10 -3 4 -2 10 17 3 0 0 4 4 14 11 2 1 0 5 0 7 5 3 6 5 3 3 3 3 8 5 2 7 9 5 7 5 -1 1 -2 1 -3 3 11 3 14 6 3 2 5 5 4 21 -1 11 3 4 3 3 6 3 4 6 8 0 1 7 7 3 8 6 1 4 6 4 1 -2 3 2 8 1 1 2 11 1 2 5 3 0 2 1 3 1 2 1 3 7 1 13 3 4 3 2 2 5 0 2 -1 9 4 1 12 3 4 -1 0 0 0 1 3 9 0 2 4 4 21 8 8 10 2 7 1 3 5 5 0 5 1 3 6 3 1 6 -3 19 4 3 6 3 5 3 10 75 3 3 1 9 3 2 28 4 2 2 3 0 2 4 2 0 4 4 4 -2 3 10 2 1 3 9 6 0 3 0 4 10 1 4 4 5 10 24 19 0 3 1 2 5 2 -3 4 5 3 9 1 3 3 5 -2 1 3 -3 3 0 12 2 1 0 2 5 5 0 -1 -1 4 5 2 11 2 5 2 1 4 2 17 3 13 3 68 6 30 5 6 3 7 14 4 2 8 4 3 1 2 1 5 12 9 2 3 1 3 0 4 0 2 2 1 4 4 5 3 3 1 1 5 3 6 3 2 21 3 0 17 5 1"
$ydsWs.Range("B3").Value = "10 25 6 6 9 16 10 16 11 9 13 4 5 6 9 15 14 8 19 18 -2 16 7 26 12 17 4 9 2 49 16 3 1 20 8 -2 12 16 33 13 3 17 18 12 3 3 14 15 0 9 23 5 7 11 6 15 6 8 -1 27 11 7 6 13 11 0 8 8 9 15 -4 5 8 19 0 25 4 12 30 13 9 7 7 12 10 16 14 12 10 16 5 19 -2 22 6 -3 16 6 6 3 1 20 21 12 16 11 8 8 3 13 17 8 15 7 16 8 28 10 15 9 9 13 5 22 9 8 8 7 9 5 8 33 11 4 5 5 6 13 15 -1 7 20 8 0 6 1 10 11 8 4 19 31 11 5 20 3 4 19 20 7 4 4 24 26 8 6 12 1 4 34 17 7 -2 20 9 0 6 13 30 2 42 15 19 9 8 12 16 4 11 14 4 6 10 9 50 17 13 7 11 7 17 8 -1 5 15 3 -5 10 9 7 3 6 5 1 15 0 23 38 0 5 25 31 5 13 6 4 30 2 7 6 6 6 -2 11 3 8 5 8 21 12 20 35 8 23 26 8 11 4 12 9 7 3 7 8 22 6 6 10 17 7 13 6 13 18 7 10 12 10 6 10 12 2 7 10 12 2 19 26 11 31 4 7 4 9 8 5 22 6 7 12 8 25 9 7 21 26 9 6 17 16 4 1 6 10 3 8 3 7 11 7 5 5 2 7 -1 2 28 19 3 24 -1 10 10 3 11 8 9 32 3 -1 1 -1 14 7 2 9 6 17 31 13 7 3 10 8 27 7 9 4 4 4 4 11 5 22 9 10 4 5 -5 10 7 4 5 7 17 1 16 15 5 16 7 3 11 4 12 15 9 4 4 21 13 5 10 12 15 1 12 3 21 30 0 9 7 7 15 11 6 9 17 23 4 5 7 16 9 3 3 21 20 1 -3 4 10 13 24 3 10 9 14 22 20 27 0 8 7 3 5 11 3 11 75 14 28 25 3 24 5 4 19 3 15 28 9 11 4 8 4 1 9 20 4 26 3 22 8 46 22 29 28 14 44 2 5 11 33 7 10 4 15 7 14 2 7 17 15 9 2 7 8 12 -1 4 9 41 3 28 3 7 18 2 12 7 11 20 3 12 26 23 4 8 5 1 10 19 17 0 15 3 26 48 11 5 5 17 6 9 9 19 12 6 9 5 18 6 10 10 6 6 15 3 19 2 10 11 4 22 20 16 13 12 38 5 19 5 3 21 10 7 20 9 11 5 41 4 7 7 12 7"
$ydsWs.Range("C3").Value = "3 12 6 13 14 -1 18 20 15 -1 9 7 3 14 5 3 18 8 4 21 9 6 9 16 6 4 9 7 7 12 54 9 4 1 9 38 14 21 9 19 8 18 18 22 14 6 14 6 6 1 4 34 0 3 24 11 13 5 2 7 10 8 12 1 7 13 8 22 17 11 2 14 21 24 8 5 5 10 6 19 45 6 6 6 1 41 7 4 18 27 17 16 35 16 8 7 23 15 14 1 14 19 20 7 8 15 16 9 13 23 6 35 18 3 11 6 12 20 12 8 12 5 41 11 21 6 5 26 5 5 11 4 13 17 50 23 15 20 8 6 19 12 15 8 17 13 21 19 6 11 31 2 3 8 9 9 12 6 7 4 14 8 6 18 10 11 3 1 13 6 4 44 2 22 3 16 6 25 21 18 9 24 4 14 4 6 7 13 22 9 6 4 4 14 20 4 19 4 8 5 16 9 6 8 17 4 9 13 6 7 -1 7 3 5 8 3 6 6 2 12 3 14 1 5 11 0 15 24 5 5 4 9 24 4 7 2 1 6 20 1 13 7 7 9 0 25 3 -2 14 16 9 7 9 7 6 2 8 6 5 2 20 3 7 8 15 3 15 3 14 5 10 5 8 4 8 6 1 10 4 7 6 17 10 6 2 0 2 4 24 6 13 4 17 50 15 11 7 10 10 18 3 9 9 20 27 8 10 1 8 10 0 20 6 4 15 6 11 18 15 12 21 53 10 9 2 5 11 12 9 8 8 17 5 9 18 23 6 36 4 0 30 12 3 -2 15 3 13 27 7 11 9 8 5 3 17 4 5 3 27 5 8 6 6 12 12 9 12 17 11 11 6 13 4 9 7 10 12 7 8 6 13 13 28 3 6 0 27 8 1 4 28 26 11 11 9 7 23 18 11 8 8 10 7 11 7 67 9 16 13 37 7 15 5 -5 27 27 27 3 14 8 11 14 11 18 17 3 23 1 13 6 13 44 20 11 3 11 4 13 0 18 27 20 1 13 12 33 3 2 4 12 24 14 6 9 7 35 12 13 7 -2 6 15 3 7 22 5 17 21 8 7 21 8 8 12 15 6 7 3 4 6 6 16 12 41 6 1 13 5 5 9 11 11 15 14 7 6 18 11 14 24 13 6 2 6 5 33 13 -3 6 21 19 17 6 9 7 12 6 13 11 11 -1 9 6 9 2 11 5 18 3 9 5 18 6 12 1 9 8 9 9 7 0 10 0 6 5 11 4 19 16 17 4 8 10 2 14 8 1 2 3 24 7 14 -3 5 8 14 16 6 11 19 7 11 7 26 5 7 6"

# ---------------------------------------------------------------
# OFF sheet - Week 13 totals roll-up (Home row 2, Road row 3)
# ---------------------------------------------------------------
$offWs = $wb.Worksheets.Item("OFF")
$offWs.Range("C2").Value = 384
$offWs.Range("E2").Value = 18
$offWs.Range("F2").Value = 117
$offWs.Range("G2").Value = 114
$offWs.Range("H2").Value = 13
$offWs.Range("I2").Value = 13
$offWs.Range("J2").Value = 77
$offWs.Range("O2").Value = 24
$offWs.Range("P2").Value = 14
$offWs.Range("C3").Value = 254
$offWs.Range("E3").Value = 60
$offWs.Range("H3").Value = 34
$offWs.Range("L3").Value = 393
$offWs.Range("M3").Value = 262
$offWs.Range("Q3").Value = 867

# ---------------------------------------------------------------
# DEF sheet - Week 13 totals roll-up (Home row 2, Road row 3)
# ---------------------------------------------------------------
$defWs = $wb.Worksheets.Item("DEF")
$defWs.Range("C2").Value = 357
$defWs.Range("E2").Value = 16
$defWs.Range("F2").Value = 108
$defWs.Range("G2").Value = 118
$defWs.Range("I2").Value = 15
$defWs.Range("J2").Value = 64
$defWs.Range("N2").Value = 28
$defWs.Range("O2").Value = 40
$defWs.Range("C3").Value = 296
$defWs.Range("E3").Value = 45
$defWs.Range("F3").Value = 175
$defWs.Range("G3").Value = 66
$defWs.Range("H3").Value = 43
$defWs.Range("I3").Value = 90
$defWs.Range("J3").Value = 86
$defWs.Range("L3").Value = 431
$defWs.Range("M3").Value = 270
$defWs.Range("Q3").Value = 811

# ---------------------------------------------------------------
# ST sheet - append Week 13 KO/PT distance entries + update totals
# ---------------------------------------------------------------
$stWs = $wb.Worksheets.Item("ST")
$stWs.Range("B2").Value = 153
$stWs.Range("D2").Value = 95
$stWs.Range("H2").Value = 6
$stWs.Range("I2").Value = 4
$stWs.Range("J2").Value = 214
$stWs.Range("K2").Value = 198
$stWs.Range("L2").Value = 119
$stWs.Range("M2").Value = 88
$stWs.Range("B3").Value = 75
$stWs.Range("B4").Value = "64 66 64 64 66 62 61 65 62 66 66 64 67 62 67 62 60 62 64 60 58 66 68 52 62 50 66 67 64 60 61 54 65 61 61 62 56 63 64 51 68 26 64 64 67 67 60 63 61 59 64 67 64 67 65 67 66 66 61 68 56 54 54 60 66 65 64 65 65 64 54 67 62 64 61 67 63 56"
$stWs.Range("B5").Value = "24 44 21 22 12 15 13 21 18 20 28 20 25 15 28 24 22 18 15 18 16 26 31 13 21 13 11 49 20 22 16 0 38 18 21 23 16 0 13 10 24 0 26 22 24 25 22 27 27 22 23 25 19 25 17 25 29 21 22 19 20 12 22 18 28 19 25 28 21 16 17 27 18 20 22 30 16 27"
$stWs.Range("D3").Value = "40 24 55 57 56 41 54 51 58 32 60 53 44 46 25 64 49 35 59 35 46 52 40 53 59 39 50 59 60 65 40 39 53 59 54 32 40 71 54 56 37 62 49 46 50 42 57 56 47 37 45 38 65 40 37 62 32 48 44 35 60 64 52 58 41 68 37 49 44 59 58 58 44 41 42 34 35 46 38 62 55 42 55 66 44 42 60 47 48 36 15 46 45 39 71"
$stWs.Range("D4").Value = "0 0 -3 0 20 0 0 0 0 0 0 0 0 0 0 0 0 0 14 0 0 0 0 7 7 0 0 7 0 0 0 0 0 8 0 0 0 0 0 0 0 7 6 0 12 0 7 0 0 0 0 0 0 0 0 18 0 6 0 0 0 25 0 0 0 0 0 10 0 5 12 11 -1 0 -1 0 0 0 0 12 0 0 4 0 0 0 11 0 0 8 0 0 0 0 0"
$stWs.Range("D5").Value = "0 0 0 0 0 12 0 11 0 0 12 0 0 0 0 0 0 0 0 15 0 0 0 0 0 0 0 5 3 0 9 70 0 14 61 0 19 3 21 6 4 3 13 0 0 5 0 17 8 12 0 7 0 -1 14 0 0 8 0 18 12 0 0 11 6 0 13 0 23 17 13 21 22 10 0 27 0 18 15 7 0 0 0 0 0 10 0 0 17 0 15 8 0 0 0 0"

# ---------------------------------------------------------------
# TURNS sheet - Week 13 totals roll-up (Road row 3)
# ---------------------------------------------------------------
$turnsWs = $wb.Worksheets.Item("TURNS")
$turnsWs.Range("D3").Value = 15
$turnsWs.Range("E3").Value = 16

# ---------------------------------------------------------------
# PEN sheet - Week 13 totals roll-up
# ---------------------------------------------------------------
$penWs = $wb.Worksheets.Item("PEN")
$penWs.Range("B2").Value = 17
$penWs.Range("D2").Value = 14
$penWs.Range("B3").Value = 29
